# "Generate Report for Handback" -- refresh the handback status report:
#   * the zh-cn / de-de handback for the second file (66fe2548-...) has now
#     come back from the localizer and its content no longer matches en-US,
#     so every "Handed back: in sync with en-US" status cell (Overview +
#     both per-locale "Status" columns) flips to "Handed back: not in sync
#     with en-US", and the status column is widened so the longer text
#     still fits.
#   * the "Correspond Handback DateTime" for that same row (row 3) is
#     stamped with the new handback timestamp on each locale sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: not in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Widen the now-longer status columns to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 32.6666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 32.6666666666667

# --- zh-cn sheet: Status column + updated handback datetime for row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K3").Value = "2016-09-09 07:12:29"
$wsZhCn.Columns.Item(3).ColumnWidth = 32.6666666666667

# --- de-de sheet: Status column + updated handback datetime for row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K3").Value = "2016-09-09 07:12:47"
$wsDeDe.Columns.Item(3).ColumnWidth = 32.6666666666667

Write-Output "Generated handback status report: updated sync status and handback datetimes"
